# New .ttl from Google sheet has been generated.
#
# The "DateType" sub-vocabulary used to have a single combined term
# "AvailableCopyrighted" at row 109. The regenerated vocabulary splits
# that into two separate terms - "Available" (row 109) and
# "Copyrighted" (row 110) - which pushes every subsequent row down by
# one. The net effect on the sheet is a single new row appended at the
# very end (old last row 162 -> new last row 163) plus the content
# edits described above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 109: this shifts the existing rows 109..162
# down to 110..163 (carrying their A/B/F inline-string content with
# them), and leaves a blank row 109 behind.
$ws.Rows.Item(109).Insert()

# Row 109 becomes the new "Available" term (split out of the old
# "AvailableCopyrighted" term), still filed under datacite:DateType.
$ws.Cells.Item(109, 1).Value = "datacite:Available"
$ws.Cells.Item(109, 2).Value = "Available"
$ws.Cells.Item(109, 6).Value = "datacite:DateType"

# Row 110 currently holds the shifted-down old row 109 content
# ("AvailableCopyrighted"/"Available Copyrighted"/"datacite:DateType").
# Overwrite it with the new "Copyrighted" term, which (per the
# regenerated sheet) no longer carries a DateType back-reference in F.
$ws.Cells.Item(110, 1).Value = "datacite:Copyrighted"
$ws.Cells.Item(110, 2).Value = "Copyrighted"
$ws.Cells.Item(110, 6).Value = ""
